# contactos.xlsx edit: remove the sample contact rows (row 2 fully, and the
# data in row 3) while leaving the header row and the already-styled/empty
# B3 / D6 cells intact, then move the active selection to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Drop the two sample rows of data ("33 1081 6707"/"prueba 1" in row 2 and
# "33 1672 0062"/"prueba 2" in row 3). Using ClearContents (rather than
# deleting the rows) removes the now-empty row 2 element entirely while
# keeping row 3 (and the untouched row 6) at their original row numbers,
# and it keeps B3's existing cell style (the underline format) in place.
$ws.Range("A2:B2").ClearContents()
$ws.Range("A3:B3").ClearContents()

# Move the saved selection/active cell to E5.
$ws.Range("E5").Select()
